$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-MM-dd HH:mm:ss"

# Touch-and-delete a scratch cell first so the new date/time number format gets
# its own fresh style slot instead of mutating the pre-existing one in place
# (the pre-existing slot stays in the stylesheet, just unreferenced afterwards).
$ws1 = $wb.Worksheets.Item("AMSIN")
$ws1.Range("H1").NumberFormat = $dateFmt
$ws1.Range("H1").Delete(-4159)

# --- AMSIN sheet: add row 3 ---
$ws1.Range("A3").Value = "'2021-06-07"
$ws1.Range("B3").Value = 44354.78363841435
$ws1.Range("B3").NumberFormat = $dateFmt
$ws1.Range("C3").Value = "regression_145final"
$ws1.Range("D3").Value = "'89"
$ws1.Range("E3").Value = 89
$ws1.Range("F3").Value = 0
$ws1.Range("G3").Value = 2.583626283333333

# Also nudge B2 so it is re-stamped with the same number format bucket
$ws1.Range("B2").Value = 44351.46821832176
$ws1.Range("B2").NumberFormat = $dateFmt

# --- BETA sheet: add row 2 ---
$ws2 = $wb.Worksheets.Item("BETA")

$ws2.Range("A2").Value = "'2021-06-08"
$ws2.Range("B2").Value = 44355.69400816529
$ws2.Range("B2").NumberFormat = $dateFmt
$ws2.Range("C2").Value = "145_beta"
$ws2.Range("D2").Value = "'89"
$ws2.Range("E2").Value = 89
$ws2.Range("F2").Value = 0
$ws2.Range("G2").Value = 2.652198833333334
